# Swap the order of "Recorded By" author names in column G:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
# Applies to every row in the used range that currently holds the old value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$updated = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
        $updated++
    }
}

Write-Host "Updated $updated cell(s) in column G."
